$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target: cell B11 changes from the text "R40" to the text "1" (still a
# plain string, not a number) while keeping its existing cell style/format
# untouched. Assigning a numeric-looking string straight to .Value makes
# Excel auto-convert it to a Number, and forcing text via a leading
# apostrophe (or via NumberFormat "@") both stamp a *new* cell style
# (quotePrefix / text format) onto the cell. So: stash the cell's current
# formatting on a scratch cell far outside the used range, write the text
# value (accepting the incidental style bump that causes), then paste the
# original formatting back on top and clean up the scratch cell.

$target = $ws.Range("B11")
$scratch = $ws.Range("Z100")

$target.Copy()
$scratch.PasteSpecial(-4122) # xlPasteFormats

$target.Value = "'1"

$scratch.Copy()
$target.PasteSpecial(-4122) # xlPasteFormats

$scratch.Clear()
